$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ESUfRaLCD-dispatch")

# Insert a new row at row 5 (shifts existing rows 5+ down by one)
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "hydro"
$ws.Range("B5").Value = "hydro es"
$ws.Range("C5").Value = "hydro dispatch"

# Move the active selection to A6, matching where the editor's cursor
# ended up after the edit (this also records A6 as this sheet's saved
# selection for next time it is opened)
$null = $ws.Range("A6").Select()

# Restore "About" as the active/visible tab, since that's what was
# showing when the workbook was saved
$null = $wb.Worksheets.Item("About").Activate()
